# Apply the edit described by the diff:
#  1. Rename the "MODEL_CONDITION" header text to "MODELCONDITION"
#  2. Delete the (now unused) first column (A), shifting B:F left to A:E

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text first (it currently lives in column E, row 1)
$ws.Range("E1").Value = "MODELCONDITION"

# Delete entire column A -- shifts columns B:F to A:E
$ws.Range("A:A").Delete()
